$d = $word.ActiveDocument

# Locate the paragraph ending with "End list arguments> ::= ')';"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*End list arguments*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $r.Collapse(0)  # wdCollapseEnd
    $r.InsertParagraphAfter()
    $newRange = $target.Next().Range
    $newRange.Text = '<case enumerator> ::= ":" ;'
}
